# Scenario3_GymsList.xlsx update
# Replaces the gym-name list in column A (rows 2-11) with a new set of
# names, and appends a new row 12, per the updated sharedStrings table.
#
# New shared strings introduced by this edit (as referenced by the sheet):
#   - D'Shoolin Functional Strength & therapy   (A2)
#   - Get the List of Top\nGyms                 (A12)
#   - Right Now Fitness                         (A9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value  = "Gyms"
$ws.Range("A2").Value  = "D'Shoolin Functional Strength & therapy"
$ws.Range("A3").Value  = "H2O Fitness Pro"
$ws.Range("A4").Value  = "Fit Lean Fitness"
$ws.Range("A5").Value  = "Ms Shine Fitness Ladies Gym"
$ws.Range("A6").Value  = "N Fit Ladies Studio (Only Ladies Gym)"
$ws.Range("A7").Value  = "Nawaz Fitness Pro"
$ws.Range("A8").Value  = "Vaishnav Fitness Zone"
$ws.Range("A9").Value  = "Right Now Fitness"
$ws.Range("A10").Value = "Fitness and Food"
$ws.Range("A11").Value = "Prime Fitness World"
$ws.Range("A12").Value = "Get the List of Top`nGyms"
